$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.21161150932312
$ws.Range("B1").Value = 2.440092086791992
$ws.Range("C1").Value = 4.737422466278076
$ws.Range("D1").Value = 2.57009220123291
$ws.Range("E1").Value = 1.086347103118896
